$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "Task Code" -> "Job Code"
$ws.Range("B3").Value = "Job Code"

# Data row: "T3" -> "03-job-cpi"
$ws.Range("B4").Value = "03-job-cpi"

# Error message text changed
$ws.Range("F4").Value = "Was not possible to read None file"

# Date column now holds a raw numeric (serial date) value instead of a text date string
$ws.Range("G4").Value = 44855.73460388806
